# Team Meeting Attendance Sep 22,2021
# Record the Sep 22 team meeting on the TEAM sheet: add the new date
# column (E) and mark who attended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TEAM")
$ws.Activate() | Out-Null

# New meeting date header (month + day), same style as the existing columns.
$ws.Range("E3").Value = "Sep"
$ws.Range("E4").Value = 22

# Mark attendance (checkmark) for everyone who showed up. Justin Murillo
# (row 9) was not in attendance, so E9 is left blank.
$ws.Range("E5").Value  = [char]0x2714
$ws.Range("E6").Value  = [char]0x2714
$ws.Range("E7").Value  = [char]0x2714
$ws.Range("E8").Value  = [char]0x2714
$ws.Range("E10").Value = [char]0x2714
$ws.Range("E11").Value = [char]0x2714

# Leave the cursor where the last entry was made.
$ws.Range("E11").Select() | Out-Null
